# Apply repulled data updates to column F (dSF) for specific rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 8
$ws.Range("F8").Value = 5
$ws.Range("F10").Value = -5
$ws.Range("F12").Value = 8
$ws.Range("F13").Value = -1
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = -7
